$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Sub Buton" -> "Sub Tümleşik Buton" (row 41, Alt Grup column E)
#    Editing the cell that is the sole owner of that shared string updates
#    the shared string text in place.
# ---------------------------------------------------------------------------
$ws.Range("E41").Value = "Sub Tümleşik Buton"

# ---------------------------------------------------------------------------
# 2) Row 32 (Alt Grup, column E) now points at the same "Sub Tümleşik Buton"
#    text instead of "4lü Tümleşik Buton".
# ---------------------------------------------------------------------------
$ws.Range("E32").Value = "Sub Tümleşik Buton"

# ---------------------------------------------------------------------------
# 3) New row 42: MSW button entry.
#    Order of the first-touched new strings matters only for matching the
#    shared-string table ordering (MSW, code, url) - not functionally
#    important, but keep it tidy / deterministic.
# ---------------------------------------------------------------------------
$ws.Range("A42").Value = 40
$ws.Range("C42").Value = "MSW"
$ws.Range("B42").Value = "BT-MSW-00-00-SLD-H3B0-03"
$ws.Range("D42").Value = "Yok"
$ws.Range("E42").Value = "Yok"
$ws.Range("F42").Value = "Tek işlevli Led"
$ws.Range("G42").Value = "Seri/Paralel"
$ws.Range("H42").Value = "Kablolu/Vidalı"
$ws.Range("I42").Value = "Buzzersız"
$ws.Range("J42").Value = "Model-03"
$ws.Range("K42").Value = "https://github.com/btk42/BT-MSW-00-00-SLD-H3B0-03"

# Hyperlink for the new Link cell (also registers the external relationship).
$ws.Hyperlinks.Add($ws.Range("K42"), "https://github.com/btk42/BT-MSW-00-00-SLD-H3B0-03")

# ---------------------------------------------------------------------------
# 4) Formatting for the new row: match row 41 (style s="11") for most
#    columns, row 2's F:H (style s="5") and K (style s="7", hyperlink look).
#    Pasting formats last also normalises the style that Hyperlinks.Add
#    applies back to the sheet's existing "Köprü"-based style (s="7").
# ---------------------------------------------------------------------------
$ws.Range("A41:E41").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)
$ws.Range("I41:J41").Copy()
$ws.Range("I42:J42").PasteSpecial(-4122)

$ws.Range("F2:H2").Copy()
$ws.Range("F42:H42").PasteSpecial(-4122)

$ws.Range("K2").Copy()
$ws.Range("K42").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) View state: zoom + selection as saved in the authored workbook.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("O30").Select()
